$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "98.715.19"
$ws.Range("E2").Value = "  +1.25%  "

# Row 3
$ws.Range("D3").Value = "3.349.24"
$ws.Range("E3").Value = "  +1.53%  "

# Row 4
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.11"
$ws.Range("E5").Value = "  +6.97%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "628.31"
$ws.Range("E6").Value = "  +1.54%  "

# Row 7
$ws.Range("E7").Value = "  +30.52%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.396"
$ws.Range("E8").Value = "  +2.98%  "

# Row 9
$ws.Range("E9").Value = "  -0.20%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.926"
$ws.Range("E10").Value = "  +15.75%  "

# Row 11
$ws.Range("D11").Value = "3.349.52"
$ws.Range("E11").Value = "  +1.68%  "

# Row 12
$ws.Range("E12").Value = "  +0.01%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.33"
$ws.Range("E13").Value = "  +12.13%  "

# Row 14
$ws.Range("D14").Value = "98.490.99"
$ws.Range("E14").Value = "  +0.65%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000251"
$ws.Range("E15").Value = "  +2.38%  "

# Row 16
$ws.Range("D16").Value = "3.975.39"
$ws.Range("E16").Value = "  +3.18%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.55"
$ws.Range("E17").Value = "  +1.07%  "

# Row 18
$ws.Range("D18").Value = "3.350.74"
$ws.Range("E18").Value = "  +1.54%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.60"
$ws.Range("E19").Value = "  -0.64%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.51"
$ws.Range("E20").Value = "  +2.78%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "489.45"
$ws.Range("E21").Value = "  -1.64%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.22"
$ws.Range("E22").Value = "  +5.27%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000209"
$ws.Range("E23").Value = "  +0.23%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.54"
$ws.Range("E24").Value = "  +4.29%  "

# Row 25
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.68"
$ws.Range("E25").Value = "  +1.30%  "

# Row 26
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "89.68"
$ws.Range("E26").Value = "  +1.07%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.03"
$ws.Range("E27").Value = "  +0.42%  "

# Row 28
$ws.Range("E28").Value = "  +24.30%  "

# Row 29
$ws.Range("D29").Value = "3.533.06"
$ws.Range("E29").Value = "  +3.46%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.19%  "

# Row 31
$ws.Range("E31").Value = "  +5.72%  "

# Row 32
$ws.Range("E32").Value = "  +12.39%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.87"
$ws.Range("E33").Value = "  +6.40%  "

# Row 34
$ws.Range("E34").Value = "  -0.08%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "28.46"
$ws.Range("E35").Value = "  +3.67%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.35"
$ws.Range("E36").Value = "  -1.40%  "

# Row 37
$ws.Range("E37").Value = "  -3.70%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.96"
$ws.Range("E38").Value = "  +1.79%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.470"
$ws.Range("E39").Value = "  +5.18%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "502.40"
$ws.Range("E40").Value = "  +2.50%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "24.90"
$ws.Range("E41").Value = "  +1.28%  "

# Row 42
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.27"
$ws.Range("E42").Value = "  +0.93%  "

# Row 43
$ws.Range("B43").Value = "MantraDAO"
$ws.Range("C43").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.69"
$ws.Range("E43").Value = "  +5.70%  "

# Row 44
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.818"
$ws.Range("E44").Value = "  +11.36%  "

# Row 45
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.35"
$ws.Range("E45").Value = "  +2.83%  "

# Row 46
$ws.Range("E46").Value = "  -0.01%  "

# Row 47
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.96"
$ws.Range("E47").Value = "  +1.13%  "

# Row 48
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "159.23"
$ws.Range("E48").Value = "  -1.43%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.863"
$ws.Range("E49").Value = "  +9.08%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.73"
$ws.Range("E50").Value = "  +3.12%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.25"
$ws.Range("E51").Value = "  +25.38%  "

